# "Alteração Na Ordem dos Slides" (Change in slide order)
#
# Swap the two adjacent slides currently at positions 14 and 15:
# the slide at position 15 ("Base de Dados - Modelagem") moves to
# position 14, and the slide that was at position 14 ("Imagem 4" /
# red-bar cover) shifts down to position 15.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$s.MoveTo(14)
